# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2404
#   *_new  -> *_FV2410
# and turn the data range into a real Excel Table, with the header row
# frozen in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) ------------------------------------
$lastCol = 21   # columns A..U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val -like '*_old') {
        $cell.Value = ($val -replace '_old$', '_FV2404')
    } elseif ($val -like '*_new') {
        $cell.Value = ($val -replace '_new$', '_FV2410')
    }
}

# --- 2. Turn the used range into an Excel Table (ListObject) -------------
$lastRow = 89
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header row renamed, table '$($tbl.Name)' created over $($dataRange.Address()), header row frozen."
